# Apply the author's edit: add two new resource rows (104 and 105) to the
# "COVID Resources-HCP" sheet, just below the existing "Health Care
# Assistants" resource block, and update the sheet selection to reflect
# where the user's cursor ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COVID Resources-HCP")

# ---------------------------------------------------------------------
# Row 104: "Well-being and Psychological Supports" (Fraser Health)
# ---------------------------------------------------------------------
$ws.Range("A104").Value = "Health Care Assistants"
$ws.Range("B104").Value = "Healthcare Provider Wellness"
$ws.Range("C104").Value = "British Columbia"
$ws.Range("D104").Value = "Well-being and Psychological Supports"
$ws.Range("E104").Value = "Fraser Health"
$ws.Range("F104").Value = "Practice Support Tool"
$ws.Hyperlinks.Add($ws.Range("G104"), "https://www.fraserhealth.ca/employees/clinical-resources/coronavirus-information/well-being-and-psychological-supports")
$ws.Range("G104").Style = "Hyperlink"
$ws.Range("G104").WrapText = $true
$ws.Rows.Item(104).RowHeight = 48

# ---------------------------------------------------------------------
# Row 105: "Home and Community Care COVID-19 Toolkit" (FNHA)
# ---------------------------------------------------------------------
$ws.Range("A105").Value = "Health Care Assistants"
$ws.Range("B105").Value = "Healthcare Provider Wellness"
$ws.Range("C105").Value = "British Columbia"
$ws.Range("D105").Value = "Home and Community Care COVID-19 Toolkit"
$ws.Range("E105").Value = "First Nations Health Authority"
$ws.Range("F105").Value = "Practice Support Tool"
$ws.Hyperlinks.Add($ws.Range("G105"), "https://www.fnha.ca/Documents/FNHA-Home-and-Community-Care-COVID-19-Toolkit.pdf")
$ws.Range("G105").Style = "Hyperlink"
$ws.Range("G105").WrapText = $true
$ws.Rows.Item(105).RowHeight = 32

# ---------------------------------------------------------------------
# Reflect the final cursor position left behind by the edit.
# ---------------------------------------------------------------------
[void]$ws.Range("C110").Select()
